# Regenerate merged AHB files:
#  - rename header columns 1-10 from "..._old" to "..._FV2304"
#  - rename header columns 12-21 from "..._new" to "..._FV2310"
#  - freeze the header row
#  - wrap the data range in an Excel Table named "Table1"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = $cell.Value()
    if ($text.EndsWith("_old")) {
        $base = $text.Substring(0, $text.Length - 4)
        $cell.Value = $base + "_FV2304"
    }
}

for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = $cell.Value()
    if ($text.EndsWith("_new")) {
        $base = $text.Substring(0, $text.Length - 4)
        $cell.Value = $base + "_FV2310"
    }
}

# Freeze panes above row 2 (keeps the header row visible)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into a native Excel table ("Table1")
$range = $ws.Range("A1:U62")
$listObject = $ws.ListObjects.Add(1, $range, 0, 1)
$listObject.Name = "Table1"
